$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.394.98"
$ws.Range("E2").Value = "  -2.40%  "
$ws.Range("D3").Value = "2.575.30"
$ws.Range("E3").Value = "  -3.01%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'545.85"
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("D6").Value = "'143.30"
$ws.Range("E6").Value = "  -2.14%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("D9").Value = "'6.79"
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("E10").Value = "  -3.92%  "
$ws.Range("E11").Value = "  +3.70%  "
$ws.Range("E12").Value = "  -2.53%  "
$ws.Range("D13").Value = "3.027.59"
$ws.Range("E13").Value = "  -3.09%  "
$ws.Range("D14").Value = "58.317.13"
$ws.Range("E14").Value = "  -2.39%  "
$ws.Range("D15").Value = "'20.48"
$ws.Range("E15").Value = "  -3.69%  "
$ws.Range("D16").Value = "2.563.65"
$ws.Range("E16").Value = "  -4.06%  "
$ws.Range("E17").Value = "  -3.99%  "
$ws.Range("D18").Value = "'4.40"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "'333.09"
$ws.Range("E19").Value = "  -3.07%  "
$ws.Range("D20").Value = "'9.97"
$ws.Range("E20").Value = "  -4.47%  "
$ws.Range("E21").Value = "  -4.85%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  -5.50%  "
$ws.Range("E27").Value = "  -4.70%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "0.0₃0733"
$ws.Range("E29").Value = "  -3.44%  "
$ws.Range("D30").Value = "'1.65"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("D31").Value = "'154.60"
$ws.Range("E31").Value = "  +2.75%  "
$ws.Range("D32").Value = "'5.84"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("E34").Value = "  -4.30%  "
$ws.Range("D35").Value = "'37.16"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "'0.840"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  -4.93%  "
$ws.Range("D38").Value = "'0.813"
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("E39").Value = "  -4.12%  "
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("D41").Value = "'277.19"
$ws.Range("E41").Value = "  -5.26%  "
$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("D45").Value = "'0.0939"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("E46").Value = "  -2.86%  "
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("D48").Value = "1.896.47"
$ws.Range("E48").Value = "  -4.32%  "
$ws.Range("E49").Value = "  -4.69%  "
$ws.Range("D50").Value = "'17.58"
$ws.Range("E50").Value = "  -4.88%  "
$ws.Range("D51").Value = "'111.24"
$ws.Range("E51").Value = "  +0.73%  "
